$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-19 17:36:05"
$ws.Range("E3").Value = "2026-02-19 17:36:07"
$ws.Range("E4").Value = "2026-02-19 17:36:10"
$ws.Range("E5").Value = "2026-02-19 17:36:12"
$ws.Range("E6").Value = "2026-02-19 17:36:15"
$ws.Range("E7").Value = "2026-02-19 17:36:17"
$ws.Range("E8").Value = "2026-02-19 17:36:20"
$ws.Range("E9").Value = "2026-02-19 17:36:22"
$ws.Range("E10").Value = "2026-02-19 17:36:25"
$ws.Range("E11").Value = "2026-02-19 17:36:27"
$ws.Range("E12").Value = "2026-02-19 17:36:30"
$ws.Range("E13").Value = "2026-02-19 17:36:32"
$ws.Range("E14").Value = "2026-02-19 17:36:35"
$ws.Range("E15").Value = "2026-02-19 17:36:37"
$ws.Range("E16").Value = "2026-02-19 17:36:39"
$ws.Range("E17").Value = "2026-02-19 17:36:42"
$ws.Range("E18").Value = "2026-02-19 17:36:45"
$ws.Range("E19").Value = "2026-02-19 17:36:47"
$ws.Range("E20").Value = "2026-02-19 17:36:50"
$ws.Range("E21").Value = "2026-02-19 17:36:52"
$ws.Range("E22").Value = "2026-02-19 17:36:54"
$ws.Range("E23").Value = "2026-02-19 17:36:56"
$ws.Range("E24").Value = "2026-02-19 17:36:59"
$ws.Range("E25").Value = "2026-02-19 17:37:01"
$ws.Range("E26").Value = "2026-02-19 17:37:04"
$ws.Range("E27").Value = "2026-02-19 17:37:06"
$ws.Range("E28").Value = "2026-02-19 17:37:09"
$ws.Range("E29").Value = "2026-02-19 17:37:11"
$ws.Range("E30").Value = "2026-02-19 17:37:14"
$ws.Range("E31").Value = "2026-02-19 17:37:17"
$ws.Range("E32").Value = "2026-02-19 17:37:19"
$ws.Range("E33").Value = "2026-02-19 17:37:21"
$ws.Range("E34").Value = "2026-02-19 17:37:24"
$ws.Range("E35").Value = "2026-02-19 17:37:26"
$ws.Range("E36").Value = "2026-02-19 17:37:29"
$ws.Range("E37").Value = "2026-02-19 17:37:31"
$ws.Range("E38").Value = "2026-02-19 17:37:34"
$ws.Range("E39").Value = "2026-02-19 17:37:36"
$ws.Range("E40").Value = "2026-02-19 17:37:39"
$ws.Range("E41").Value = "2026-02-19 17:37:41"
$ws.Range("E42").Value = "2026-02-19 17:37:43"
$ws.Range("E43").Value = "2026-02-19 17:37:46"
$ws.Range("E44").Value = "2026-02-19 17:37:48"
$ws.Range("E45").Value = "2026-02-19 17:37:50"
$ws.Range("E46").Value = "2026-02-19 17:37:53"
